# Updated cryptos list on Sun Oct  1 13:42:16 UTC 2023 with GitHub Actions
# Refreshes the Price (D) and Volume(1h) (E) columns of the crypto table
# on the active sheet with the latest scraped figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.206.55"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.93%  "
$ws.Range("D3").Value = "'1.687.10"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.64%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'215.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.35%  "
$ws.Range("E6").Value = "  +0.58%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").Value = "'23.23"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +10.86%  "
$ws.Range("E9").Value = "  +4.83%  "
$ws.Range("E10").Value = "  +1.43%  "
$ws.Range("D11").Value = "'0.0890"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.22%  "
$ws.Range("D12").Value = "'1.925.24"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.63%  "
$ws.Range("D13").Value = "'1.687.51"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.04%  "
$ws.Range("E14").Value = "  +2.51%  "
$ws.Range("D15").Value = "'0.554"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.51%  "
$ws.Range("D16").Value = "'67.19"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.24%  "
$ws.Range("D17").Value = "'27.204.54"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.87%  "
$ws.Range("D18").Value = "'236.15"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.10%  "
$ws.Range("E19").Value = "  -1.46%  "
$ws.Range("D20").Value = "'0.0₃0745"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.56%  "
$ws.Range("E21").Value = "  +0.10%  "
$ws.Range("D22").Value = "'4.57"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.09%  "
$ws.Range("D23").Value = "'9.67"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.18%  "
$ws.Range("D24").Value = "'2.09"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.70%  "
$ws.Range("D25").Value = "'147.25"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.25%  "
$ws.Range("E26").Value = "  +1.22%  "
$ws.Range("E27").Value = "  +2.53%  "
$ws.Range("D28").Value = "'0.113"
$ws.Range("D28").Style = "Normal"
$ws.Range("E29").Value = "  -0.10%  "
$ws.Range("E30").Value = "  +0.70%  "
$ws.Range("E31").Value = "  +0.42%  "
$ws.Range("E32").Value = "  +1.98%  "
$ws.Range("D33").Value = "'1.552.35"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.30%  "
$ws.Range("D34").Value = "'3.25"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.37%  "
$ws.Range("E35").Value = "  -0.75%  "
$ws.Range("D36").Value = "'0.949"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.35%  "
$ws.Range("E37").Value = "  +3.43%  "
$ws.Range("E38").Value = "  -0.31%  "
$ws.Range("E39").Value = "  -0.17%  "
$ws.Range("E40").Value = "  +3.30%  "
$ws.Range("D41").Value = "'69.23"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.75%  "
$ws.Range("E42").Value = "  +0.04%  "
$ws.Range("E43").Value = "  +0.09%  "
$ws.Range("D44").Value = "'2.27"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.35%  "
$ws.Range("E45").Value = "  +0.87%  "
$ws.Range("D46").Value = "'0.789"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.23%  "
$ws.Range("D47").Value = "'91.10"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.71%  "
$ws.Range("E48").Value = "  +5.83%  "
$ws.Range("E49").Value = "  +3.40%  "
$ws.Range("E50").Value = "  +7.52%  "
$ws.Range("E51").Value = "  +1.69%  "
